$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.968.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.387.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.45"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.98%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.95"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.389.83"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.39"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.20%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.40%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.975.88"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.33%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000188"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.60%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.06"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.064.42"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.339.93"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.24"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.04"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.16"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.01"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.51%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.27"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.528"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000117"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +23.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.42"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.179"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.21%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.08"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +8.58%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.55%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.12"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.37"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.57%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.76"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.18%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.86"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.34%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.956.37"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0754"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.89"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.62%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.85%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.58"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.32%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.756"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.11%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.31"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.95%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.27"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.06%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.02%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +20.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.34"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.829"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.02%  "
